# Auto update stocks_data.xlsx [2025-12-15 00:55:16]
# Appends a new column Q that mirrors column P (same trading-date block),
# duplicating most values and updating a handful of index-close cells that
# were revised after the original column P was written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the whole P column (header rows + data rows) into Q,
#    carrying over values, types and per-row styles in one shot.
$ws.Range("P1:P118").Copy($ws.Range("Q1:Q118"))

# 2) A few index-value rows were refreshed with a slightly different
#    close value than column P (column P itself is left untouched).
$ws.Cells.Item(4, 17).Value = 3884.54
$ws.Cells.Item(7, 17).Value = 5528.82
$ws.Cells.Item(10, 17).Value = 4579.49
$ws.Cells.Item(13, 17).Value = 7169.64
$ws.Cells.Item(28, 17).Value = 50254.99
$ws.Cells.Item(58, 17).Value = 15415.3
$ws.Cells.Item(100, 17).Value = 9065.18
$ws.Cells.Item(118, 17).Value = 2958.04

# 3) Two stocks had no data for this refresh even though column P had a
#    value, so their column-Q cell stays blank (format-only copy already
#    happened above; just clear the carried-over value).
$ws.Cells.Item(37, 17).Value = ""
$ws.Cells.Item(40, 17).Value = ""

# 4) Match the new column's stored width (same 15-character width as the
#    rest of the sheet) so a <col> entry for column 17 is emitted.
$ws.Columns.Item(17).ColumnWidth = 14.29
